# Workbook: Hortaliza, Terminal Hortofrutícola Agro Chillán - Perejil.xlsx
# Commit: "Fruta / hortaliza, semanal"
#
# The diff shows a new weekly record inserted as row 7 (pushing the
# previous rows 7..117 down to 8..118, i.e. a classic "insert row"
# operation), with the new row populated with a fresh observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 7; this shifts existing rows 7-117 down to 8-118
# and copies formatting (incl. the date style used in column D) from the
# row above, same as Excel's native "Insert Copied/Sheet Rows" behavior.
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new weekly observation.
$ws.Range("A7").Value = 7
$ws.Range("B7").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C7").Value = "Ñuble"
$ws.Range("D7").Value = 45190
$ws.Range("E7").Value = 16
$ws.Range("F7").Value = 100112044
$ws.Range("G7").Value = "Perejil"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 150
$ws.Range("K7").Value = 1500
$ws.Range("L7").Value = 1500
$ws.Range("M7").Value = 1500
$ws.Range("N7").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O7").Value = "Región de Ñuble"
$ws.Range("P7").Value = 1500
$ws.Range("Q7").Value = 1
$ws.Range("R7").Value = "Hortaliza"
